$d = $word.ActiveDocument

$pairs = @(
    @{old="498÷3="; new="919÷6="},
    @{old="612÷3="; new="969÷2="},
    @{old="140÷8="; new="855÷2="},
    @{old="976÷8="; new="717÷6="},
    @{old="220÷6="; new="421÷3="},
    @{old="197÷8="; new="397÷4="},
    @{old="805÷5="; new="653÷8="},
    @{old="609÷7="; new="976÷4="},
    @{old="818÷8="; new="894÷5="},
    @{old="516÷8="; new="864÷8="},
    @{old="487÷3="; new="667÷2="},
    @{old="877÷2="; new="171÷5="},
    @{old="111÷9="; new="535÷2="},
    @{old="304÷4="; new="822÷9="},
    @{old="965÷2="; new="760÷6="},
    @{old="378÷8="; new="873÷8="},
    @{old="650÷4="; new="290÷4="},
    @{old="800÷6="; new="398÷5="},
    @{old="562÷6="; new="588÷2="},
    @{old="489÷9="; new="642÷5="},
    @{old="504÷9="; new="717÷4="},
    @{old="876÷9="; new="231÷3="},
    @{old="927÷9="; new="155÷3="},
    @{old="564÷2="; new="899÷9="},
    @{old="581÷9="; new="551÷7="}
)

foreach ($pair in $pairs) {
    $rng = $d.Content
    $rng.Find.Execute($pair.old, $true, $false, $false, $false, $false, $true, 1, $false, $pair.new, 2)
}
